# Add an "Address" column, inserted before the existing "District" column (F),
# pushing District from F to G. The new Address column is populated with the
# school/institution address extracted from each teacher's row (column B),
# matching the target workbook produced by the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; this shifts the existing column F ("District")
# and its data one column to the right, becoming column G.
$ws.Columns("F:F").Insert()

# New header for the inserted column.
$ws.Range("F2").Value = "Address"

# Per-row Address values (row number -> text). Rows not present here are
# left blank, matching rows where the source address text did not
# cleanly resolve to a single address fragment.
$addr = @{}
$addr[3]  = "G H S KodiyalaGubbi"
$addr[4]  = "Jnanapeeta High School Chikkanayakanahalli"
$addr[5]  = "S M H S Kora"
$addr[6]  = "S S P U C High School SectionB H Road"
$addr[7]  = "G H S Kyathsandra"
$addr[8]  = "S R H S KonehallyTiptur"
$addr[9]  = "K N D H S HuliyarChikkanayakanahalli"
$addr[10] = "Govt. High School AlburTiptur"
$addr[11] = "G H S HulikalTuruvekere"
$addr[12] = "G H S NeralekereGubbi"
$addr[13] = "Sri Bhyraveswara RuralHigh School Bugudanahalli"
$addr[15] = "Sree Siddaganga Rural High School Halanur"
$addr[16] = "R S M T High School"
$addr[17] = "G H S HonnashettyhallyGubbi"
$addr[18] = "G P U C(H S) Huliyar – KenkereChikkanayakanahalli"
$addr[19] = "G H S ThyagaturGubbi"
$addr[20] = "G H S Yalagalavadi"
$addr[22] = "Sri Someshwara High School Dombarana HalliTuruvekere"
$addr[24] = "Sree Siddaganga Rural High School"
$addr[25] = "G J C Biligere Tiptur"
$addr[26] = "S R S H S BelaguliChikkanayakanahalli"
$addr[27] = "G H S Sampige Hosahally Turuvekere"
$addr[28] = "G J C KadabaGubbi"
$addr[29] = "Government High SchoolSanthemavathurKunigal"
$addr[30] = "G H S AjjanahalliTuruvekere"
$addr[31] = "Siddaganga High (SSPUC) School B H Road"
$addr[32] = "G H S KodagihalliTuruvekere"
$addr[33] = "B R H S Beladhara"
$addr[34] = "Sri Adarsha High SchoolKonthihally"
$addr[35] = "G J C Dandinashivara Turvekere"
$addr[36] = "Govt. High SchoolAlkereKunigal"
$addr[37] = "G H S Halepalya Tiptur"
$addr[39] = "Turvekere"
$addr[40] = "G J C Nedasale"
$addr[41] = "G J C (HS)NagasandraKunigal"
$addr[42] = "Sree Swarnamba High School"
$addr[43] = "M G G P U C (High School Section) Kunigal"
$addr[44] = "Govt. High School BajaguruTiptur"
$addr[45] = "Siddaganga Composite PU College (High School Section) B H Road"
$addr[46] = "Govt. High School MuddanahalliTuruvekere"
$addr[47] = "G J C Gubbi"
$addr[48] = "G H S K Kallahalli Gubbi"
$addr[49] = "G H S TadasurTiptur"
$addr[50] = "Government High SchoolAremallenahalliTuruvekere"
$addr[51] = "G G J C Tiptur"
$addr[52] = "Govt. Composite PU College Oorukere"
$addr[53] = "G H S KondliGubbi"
$addr[54] = "G H S MadenahalliGubbi"
$addr[55] = "Nehru Vidya Shalam Mayasandra T B Turuvekere"
$addr[56] = "G J C KadabaGubbi"
$addr[58] = "Tiptur"
$addr[60] = "G H S ChottanahalliKunigal"
$addr[61] = "G H S Shivara Tiptur"
$addr[62] = "Vivekananda High SchoolVinoba Nagar"
$addr[63] = "Govt. High School ValagerepuraKunigal"

foreach ($r in $addr.Keys) {
    $ws.Cells.Item($r, 6).Value = $addr[$r]
}
